# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns to the
# latest scraped readings (coinranking.com feed refresh via GitHub Actions).
#
# Every cell in these two columns is stored as literal text in the workbook
# (prices use dotted thousands separators like "28.259.07", percentages keep
# padding spaces like "  -0.81%  "), so values are written back as text too.
# A handful of the new price strings parse as plain numbers (e.g. "2.00"),
# and Excel auto-converts a bare numeric-looking Value into a real number
# (dropping the trailing zero). To keep those as text we briefly use the
# leading-apostrophe text-entry convention and then restore the default
# "Normal" style so the quote-prefix formatting does not linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.259.07"   # was "28.273.39"
$ws.Range("E2").Value = "  -0.81%  "   # was "  -0.60%  "
$ws.Range("D3").Value = "1.550.54"   # was "1.551.78"
$ws.Range("E3").Value = "  -1.25%  "   # was "  -1.11%  "
$ws.Range("E4").Value = "  +0.18%  "   # was "  +0.17%  "
$ws.Range("D5").Value = "'208.87"   # was "209.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "   # was "  -1.26%  "
$ws.Range("D6").Value = "'0.485"   # was "0.486"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.51%  "   # was "  -1.33%  "
$ws.Range("E7").Value = "  +0.16%  "   # was "  +0.13%  "
$ws.Range("D8").Value = "'23.38"   # was "23.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.67%  "   # was "  -2.15%  "
$ws.Range("E9").Value = "  -2.05%  "   # was "  -1.77%  "
$ws.Range("E10").Value = "  -1.36%  "   # was "  -1.16%  "
$ws.Range("E11").Value = "  +0.15%  "   # was "  +0.12%  "
$ws.Range("D12").Value = "1.773.44"   # was "1.774.32"
$ws.Range("E12").Value = "  -1.16%  "   # was "  -1.08%  "
$ws.Range("D13").Value = "1.555.62"   # was "1.551.65"
$ws.Range("E13").Value = "  -0.88%  "   # was "  -1.05%  "
$ws.Range("D14").Value = "28.268.13"   # was "28.290.31"
$ws.Range("E14").Value = "  -0.64%  "   # was "  -0.46%  "
$ws.Range("E15").Value = "  -1.19%  "   # was "  -0.96%  "
$ws.Range("E16").Value = "  -2.41%  "   # was "  -2.21%  "
$ws.Range("D17").Value = "'60.22"   # was "60.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.03%  "   # was "  -2.91%  "
$ws.Range("D18").Value = "'228.35"   # was "227.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.79%  "   # was "  -0.53%  "
$ws.Range("D19").Value = "'7.28"   # was "7.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "   # was "  -0.33%  "
$ws.Range("E20").Value = "  -2.81%  "   # was "  -2.47%  "
$ws.Range("E21").Value = "  +0.15%  "   # was "  +0.18%  "
$ws.Range("E22").Value = "  +0.58%  "   # was "  +0.92%  "
$ws.Range("E23").Value = "  -3.18%  "   # was "  -2.95%  "
$ws.Range("D24").Value = "'2.00"   # was "2.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.29%  "   # was "  -4.78%  "
$ws.Range("D25").Value = "'147.44"   # was "147.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.14%  "   # was "  -2.10%  "
$ws.Range("E26").Value = "  -1.66%  "   # was "  -1.45%  "
$ws.Range("E27").Value = "  -0.19%  "   # was "  +0.01%  "
$ws.Range("E28").Value = "  +0.16%  "   # was "  +0.13%  "
$ws.Range("D29").Value = "'6.21"   # was "6.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.27%  "   # was "  -3.09%  "
$ws.Range("D30").Value = "'0.0466"   # was "0.0467"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.36%  "   # was "  -2.14%  "
$ws.Range("E31").Value = "  -3.96%  "   # was "  -4.07%  "
$ws.Range("D32").Value = "'3.15"   # was "3.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "   # was "  -0.77%  "
$ws.Range("E33").Value = "  -1.27%  "   # was "  -1.09%  "
$ws.Range("D34").Value = "1.381.86"   # was "1.382.61"
$ws.Range("E34").Value = "  -0.64%  "   # was "  -0.61%  "
$ws.Range("E35").Value = "  +1.07%  "   # was "  +1.03%  "
$ws.Range("E36").Value = "  -3.06%  "   # was "  -3.22%  "
$ws.Range("E37").Value = "  -0.85%  "   # was "  -0.89%  "
$ws.Range("D38").Value = "'2.59"   # was "2.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "   # was "  -1.51%  "
$ws.Range("E39").Value = "  -2.88%  "   # was "  -2.49%  "
$ws.Range("D40").Value = "'0.509"   # was "0.510"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "   # was "  -2.33%  "
$ws.Range("E41").Value = "  +1.39%  "   # was "  +1.61%  "
$ws.Range("E42").Value = "  +0.12%  "   # was "  +0.15%  "
$ws.Range("D43").Value = "'0.771"   # was "0.772"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.99%  "   # was "  -1.78%  "
$ws.Range("D44").Value = "'0.0463"   # was "0.0464"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "   # was "  -0.91%  "
$ws.Range("D45").Value = "'5.37"   # was "5.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "   # was "  -1.25%  "
$ws.Range("D46").Value = "'61.48"   # was "61.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "   # was "  -1.71%  "
$ws.Range("D47").Value = "'0.909"   # was "0.910"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.31%  "   # was "  -6.18%  "
$ws.Range("D48").Value = "1.686.61"   # was "1.687.00"
$ws.Range("D49").Value = "'85.14"   # was "85.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "   # was "  -0.95%  "
$ws.Range("E50").Value = "  +0.27%  "   # was "  +0.68%  "
$ws.Range("D51").Value = "'41.34"   # was "41.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.53%  "   # was "  +8.05%  "
